$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 24.50527833333333
$ws.Cells.Item(2, 8).Value = 73.515835
$ws.Cells.Item(2, 9).Value = 0.04846830138877924
$ws.Cells.Item(2, 10).Value = 0.04846830138877924
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 34.53682066666666
$ws.Cells.Item(2, 14).Value = 103.610462
$ws.Cells.Item(2, 15).Value = 0.2803141013583512
$ws.Cells.Item(2, 16).Value = 0.2803141013583513
$ws.Cells.Item(2, 17).Value = 846.3344031850854
$ws.Cells.Item(2, 18).Value = 7617.009628665769
$ws.Cells.Item(2, 19).Value = 0.01358634834816138
$ws.Cells.Item(2, 20).Value = 0.01358634834816138

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 24.50527833333333
$ws.Cells.Item(3, 8).Value = 73.515835
$ws.Cells.Item(3, 9).Value = 0.04846830138877924
$ws.Cells.Item(3, 10).Value = 0.04846830138877924
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 47.750315
$ws.Cells.Item(3, 14).Value = 143.250945
$ws.Cells.Item(3, 15).Value = 0.3875598963781245
$ws.Cells.Item(3, 16).Value = 0.3875598963781245
$ws.Cells.Item(3, 17).Value = 1170.134759579342
$ws.Cells.Item(3, 18).Value = 10531.21283621407
$ws.Cells.Item(3, 19).Value = 0.01878436986385899
$ws.Cells.Item(3, 20).Value = 0.01878436986385899

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 24.50527833333333
$ws.Cells.Item(4, 8).Value = 73.515835
$ws.Cells.Item(4, 9).Value = 0.04846830138877924
$ws.Cells.Item(4, 10).Value = 0.04846830138877924
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 28.69151433333333
$ws.Cells.Item(4, 14).Value = 86.07454300000001
$ws.Cells.Item(4, 15).Value = 0.232871350104353
$ws.Cells.Item(4, 16).Value = 0.232871350104353
$ws.Cells.Item(4, 17).Value = 703.0935445431561
$ws.Cells.Item(4, 18).Value = 6327.841900888405
$ws.Cells.Item(4, 19).Value = 0.01128687878166971
$ws.Cells.Item(4, 20).Value = 0.01128687878166971

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 24.50527833333333
$ws.Cells.Item(5, 8).Value = 73.515835
$ws.Cells.Item(5, 9).Value = 0.04846830138877924
$ws.Cells.Item(5, 10).Value = 0.04846830138877924
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 12.228925
$ws.Cells.Item(5, 14).Value = 36.686775
$ws.Cells.Item(5, 15).Value = 0.09925465215917123
$ws.Cells.Item(5, 16).Value = 0.09925465215917123
$ws.Cells.Item(5, 17).Value = 299.6732108424584
$ws.Cells.Item(5, 18).Value = 2697.058897582125
$ws.Cells.Item(5, 19).Value = 0.004810704395089159
$ws.Cells.Item(5, 20).Value = 0.004810704395089159

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 110.1980973333333
$ws.Cells.Item(6, 8).Value = 330.594292
$ws.Cells.Item(6, 9).Value = 0.2179577200213544
$ws.Cells.Item(6, 10).Value = 0.2179577200213544
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 34.53682066666666
$ws.Cells.Item(6, 14).Value = 103.610462
$ws.Cells.Item(6, 15).Value = 0.2803141013583512
$ws.Cells.Item(6, 16).Value = 0.2803141013583513
$ws.Cells.Item(6, 17).Value = 3805.891925409211
$ws.Cells.Item(6, 18).Value = 34253.0273286829
$ws.Cells.Item(6, 19).Value = 0.06109662242190109
$ws.Cells.Item(6, 20).Value = 0.0610966224219011

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 110.1980973333333
$ws.Cells.Item(7, 8).Value = 330.594292
$ws.Cells.Item(7, 9).Value = 0.2179577200213544
$ws.Cells.Item(7, 10).Value = 0.2179577200213544
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 47.750315
$ws.Cells.Item(7, 14).Value = 143.250945
$ws.Cells.Item(7, 15).Value = 0.3875598963781245
$ws.Cells.Item(7, 16).Value = 0.3875598963781245
$ws.Cells.Item(7, 17).Value = 5261.993860067327
$ws.Cells.Item(7, 18).Value = 47357.94474060594
$ws.Cells.Item(7, 19).Value = 0.0844716713862884
$ws.Cells.Item(7, 20).Value = 0.0844716713862884

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 110.1980973333333
$ws.Cells.Item(8, 8).Value = 330.594292
$ws.Cells.Item(8, 9).Value = 0.2179577200213544
$ws.Cells.Item(8, 10).Value = 0.2179577200213544
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 28.69151433333333
$ws.Cells.Item(8, 14).Value = 86.07454300000001
$ws.Cells.Item(8, 15).Value = 0.232871350104353
$ws.Cells.Item(8, 16).Value = 0.232871350104353
$ws.Cells.Item(8, 17).Value = 3161.750289145395
$ws.Cells.Item(8, 18).Value = 28455.75260230856
$ws.Cells.Item(8, 19).Value = 0.05075610852703939
$ws.Cells.Item(8, 20).Value = 0.05075610852703939

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 110.1980973333333
$ws.Cells.Item(9, 8).Value = 330.594292
$ws.Cells.Item(9, 9).Value = 0.2179577200213544
$ws.Cells.Item(9, 10).Value = 0.2179577200213544
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 12.228925
$ws.Cells.Item(9, 14).Value = 36.686775
$ws.Cells.Item(9, 15).Value = 0.09925465215917123
$ws.Cells.Item(9, 16).Value = 0.09925465215917123
$ws.Cells.Item(9, 17).Value = 1347.604267432034
$ws.Cells.Item(9, 18).Value = 12128.4384068883
$ws.Cells.Item(9, 19).Value = 0.02163331768612557
$ws.Cells.Item(9, 20).Value = 0.02163331768612557

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 351.7202226666666
$ws.Cells.Item(10, 8).Value = 1055.160668
$ws.Cells.Item(10, 9).Value = 0.6956575446665283
$ws.Cells.Item(10, 10).Value = 0.6956575446665284
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 34.53682066666666
$ws.Cells.Item(10, 14).Value = 103.610462
$ws.Cells.Item(10, 15).Value = 0.2803141013583512
$ws.Cells.Item(10, 16).Value = 0.2803141013583513
$ws.Cells.Item(10, 17).Value = 12147.29825507873
$ws.Cells.Item(10, 18).Value = 109325.6842957086
$ws.Cells.Item(10, 19).Value = 0.1950026194863549
$ws.Cells.Item(10, 20).Value = 0.195002619486355

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 351.7202226666666
$ws.Cells.Item(11, 8).Value = 1055.160668
$ws.Cells.Item(11, 9).Value = 0.6956575446665283
$ws.Cells.Item(11, 10).Value = 0.6956575446665284
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 47.750315
$ws.Cells.Item(11, 14).Value = 143.250945
$ws.Cells.Item(11, 15).Value = 0.3875598963781245
$ws.Cells.Item(11, 16).Value = 0.3875598963781245
$ws.Cells.Item(11, 17).Value = 16794.75142420347
$ws.Cells.Item(11, 18).Value = 151152.7628178313
$ws.Cells.Item(11, 19).Value = 0.2696089659256202
$ws.Cells.Item(11, 20).Value = 0.2696089659256203

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 351.7202226666666
$ws.Cells.Item(12, 8).Value = 1055.160668
$ws.Cells.Item(12, 9).Value = 0.6956575446665283
$ws.Cells.Item(12, 10).Value = 0.6956575446665284
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 28.69151433333333
$ws.Cells.Item(12, 14).Value = 86.07454300000001
$ws.Cells.Item(12, 15).Value = 0.232871350104353
$ws.Cells.Item(12, 16).Value = 0.232871350104353
$ws.Cells.Item(12, 17).Value = 10091.38580996386
$ws.Cells.Item(12, 18).Value = 90822.47228967473
$ws.Cells.Item(12, 19).Value = 0.1619987116367737
$ws.Cells.Item(12, 20).Value = 0.1619987116367737

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 351.7202226666666
$ws.Cells.Item(13, 8).Value = 1055.160668
$ws.Cells.Item(13, 9).Value = 0.6956575446665283
$ws.Cells.Item(13, 10).Value = 0.6956575446665284
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 12.228925
$ws.Cells.Item(13, 14).Value = 36.686775
$ws.Cells.Item(13, 15).Value = 0.09925465215917123
$ws.Cells.Item(13, 16).Value = 0.09925465215917123
$ws.Cells.Item(13, 17).Value = 4301.160223973967
$ws.Cells.Item(13, 18).Value = 38710.4420157657
$ws.Cells.Item(13, 19).Value = 0.06904724761777939
$ws.Cells.Item(13, 20).Value = 0.0690472476177794

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 19.17031833333333
$ws.Cells.Item(14, 8).Value = 57.510955
$ws.Cells.Item(14, 9).Value = 0.03791643392333802
$ws.Cells.Item(14, 10).Value = 0.03791643392333802
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 34.53682066666666
$ws.Cells.Item(14, 14).Value = 103.610462
$ws.Cells.Item(14, 15).Value = 0.2803141013583512
$ws.Cells.Item(14, 16).Value = 0.2803141013583513
$ws.Cells.Item(14, 17).Value = 662.0818464012456
$ws.Cells.Item(14, 18).Value = 5958.73661761121
$ws.Cells.Item(14, 19).Value = 0.0106285111019338
$ws.Cells.Item(14, 20).Value = 0.0106285111019338

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 19.17031833333333
$ws.Cells.Item(15, 8).Value = 57.510955
$ws.Cells.Item(15, 9).Value = 0.03791643392333802
$ws.Cells.Item(15, 10).Value = 0.03791643392333802
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 47.750315
$ws.Cells.Item(15, 14).Value = 143.250945
$ws.Cells.Item(15, 15).Value = 0.3875598963781245
$ws.Cells.Item(15, 16).Value = 0.3875598963781245
$ws.Cells.Item(15, 17).Value = 915.3887390669418
$ws.Cells.Item(15, 18).Value = 8238.498651602475
$ws.Cells.Item(15, 19).Value = 0.01469488920235689
$ws.Cells.Item(15, 20).Value = 0.01469488920235689

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 19.17031833333333
$ws.Cells.Item(16, 8).Value = 57.510955
$ws.Cells.Item(16, 9).Value = 0.03791643392333802
$ws.Cells.Item(16, 10).Value = 0.03791643392333802
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 28.69151433333333
$ws.Cells.Item(16, 14).Value = 86.07454300000001
$ws.Cells.Item(16, 15).Value = 0.232871350104353
$ws.Cells.Item(16, 16).Value = 0.232871350104353
$ws.Cells.Item(16, 17).Value = 550.0254632353962
$ws.Cells.Item(16, 18).Value = 4950.229169118566
$ws.Cells.Item(16, 19).Value = 0.008829651158870215
$ws.Cells.Item(16, 20).Value = 0.008829651158870215

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 19.17031833333333
$ws.Cells.Item(17, 8).Value = 57.510955
$ws.Cells.Item(17, 9).Value = 0.03791643392333802
$ws.Cells.Item(17, 10).Value = 0.03791643392333802
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 12.228925
$ws.Cells.Item(17, 14).Value = 36.686775
$ws.Cells.Item(17, 15).Value = 0.09925465215917123
$ws.Cells.Item(17, 16).Value = 0.09925465215917123
$ws.Cells.Item(17, 17).Value = 234.4323851244584
$ws.Cells.Item(17, 18).Value = 2109.891466120125
$ws.Cells.Item(17, 19).Value = 0.003763382460177115
$ws.Cells.Item(17, 20).Value = 0.003763382460177115

